# Implemented CR 1. And update audit file to be create monthly
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new settings row (above the old row 15 / "TransactionType") for the
# new NikeCODMasterCredentialsAsset / NikeCODMasterPassword asset pair.
$ws.Rows("15:15").Insert()
$ws.Range("A15").Value = "NikeCODMasterCredentialsAsset"
$ws.Range("B15").Value = "NikeCODMasterPassword"

# TransactionAuditFile location now points at a monthly-generated audit file.
$ws.Range("B8").Value = "C:\ESW\Audit\Nike Daily\{year}\{month} Audit NikeDailyRefund Transactions.xlsx"

# FromDate default value updated (row shifted down to 22 after the insert above).
$ws.Range("B22").Value = "01/01/2020"

# Update the saved selection to match the author's last cursor position.
$ws.Activate() | Out-Null
$ws.Range("B20").Select() | Out-Null
